# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (totals) sheet.
#    It is created by duplicating the "2021-Q4" sheet (so it inherits the
#    exact same header/index-column styling used by every other quarterly
#    sheet) and then replacing its data with the single 2022-Q1 fund
#    holding row.
# 2. Prepend a new row to the "总计" sheet summarising the new quarter
#    (date/count/value), pushing the previously existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell while forcing it to be stored as
# TEXT (inline/shared string), matching how the source data keeps
# numeric-looking strings ("18.02", "010695", ...) as text rather than
# letting Excel auto-convert them to numbers.
# ---------------------------------------------------------------------
function Set-TextValue($rng, $value) {
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# =======================================================================
# 1. Add the "2022-Q1" worksheet immediately before "总计", cloned from
#    "2021-Q4" so the cell styles match the existing quarterly sheets.
# =======================================================================
$sourceSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sourceSheet.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count - 1)
$newSheet.Name = "2022-Q1"

# "2021-Q4" had three fund-holding rows (rows 2-4); 2022-Q1 only needs one.
$newSheet.Rows.Item(3).Resize(2).Delete()

# Update the single remaining data row with the 2022-Q1 figures.
Set-TextValue $newSheet.Range("B2") "010695"
Set-TextValue $newSheet.Range("C2") "华夏磐益一年定期开放混合"
Set-TextValue $newSheet.Range("D2") "18.02"
Set-TextValue $newSheet.Range("E2") "82.41"
Set-TextValue $newSheet.Range("F2") "2.53"
Set-TextValue $newSheet.Range("G2") "0.4559"
$newSheet.Range("H2").Value = 8

# =======================================================================
# 2. Prepend a summary row for 2022-Q1 on the "总计" sheet
# =======================================================================
$totalWs = $wb.Worksheets.Item("总计")
$totalWs.Rows.Item(2).Insert()
$totalWs.Rows.Item(2).ClearFormats()

# Re-use the bold/bordered index-column style already applied to the rows
# below (e.g. the row that used to be row 2, now pushed down to row 3) by
# copying its formatting onto the new A2 cell.
$totalWs.Range("A3").Copy()
$totalWs.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("C2").Value = 1
$totalWs.Range("D2").Value = 0.46
